$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3521
$ws1.Range("F4").Value = 143
$ws1.Range("F5").Value = 7018
$ws1.Range("F6").Value = 3133
$ws1.Range("F7").Value = 55
$ws1.Range("F8").Value = 133
$ws1.Range("F9").Value = 30
$ws1.Range("F10").Value = 39
$ws1.Range("F11").Value = 87
$ws1.Range("F12").Value = 45
$ws1.Range("F13").Value = 14
$ws1.Range("F14").Value = 179
$ws1.Range("F15").Value = 588
$ws1.Range("F16").Value = 20
$ws1.Range("F17").Value = 45

# Sheet "全部类型" (sheet4) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3521
$ws4.Range("F5").Value = 143
$ws4.Range("F6").Value = 7018
$ws4.Range("F7").Value = 3133
$ws4.Range("F8").Value = 55
$ws4.Range("F9").Value = 133
$ws4.Range("F10").Value = 30
$ws4.Range("F11").Value = 39
$ws4.Range("F12").Value = 87
$ws4.Range("F13").Value = 45
$ws4.Range("F14").Value = 14
$ws4.Range("F15").Value = 179
$ws4.Range("F16").Value = 588
$ws4.Range("F17").Value = 20
$ws4.Range("F18").Value = 45
